# -----------------------------------------------------------------------
# Commit: "Consolidate ... phases" / "Consolidate ... components" /
#         "Add Presentation Notes section header" family of fixes.
#
# Concretely, on this deck it:
#   1) Un-bolds the header row of the "Timeline & Milestones" table
#      (slide "Timeline & Milestones") and the bold "Phase 3" data row.
#   2) Un-bolds the header row of the "Investment Summary" table.
#   3) On the "Next Steps" slide, bolds just the "Label:" prefix of each
#      bullet line (splitting the single run into a bold run + a plain
#      run for the remaining text).
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

function Get-TableShape($slide) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            return $shp
        }
    }
    return $null
}

function Get-SlideByTitle($pres, $title) {
    for ($i = 1; $i -le $pres.Slides.Count; $i++) {
        $slide = $pres.Slides.Item($i)
        for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
            $shp = $slide.Shapes.Item($j)
            if ($shp.HasTextFrame) {
                if ($shp.TextFrame.TextRange.Text -eq $title) {
                    return $slide
                }
            }
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1) Timeline & Milestones table: un-bold header row + "Phase 3" row
# ---------------------------------------------------------------
$slideTimeline = Get-SlideByTitle $p "Timeline & Milestones"
$tblTimeline = Get-TableShape $slideTimeline
$tbl1 = $tblTimeline.Table

for ($r = 1; $r -le $tbl1.Rows.Count; $r++) {
    $firstCellText = $tbl1.Cell($r, 1).Shape.TextFrame.TextRange.Text
    if (($r -eq 1) -or ($firstCellText -eq "Phase 3")) {
        for ($c = 1; $c -le $tbl1.Columns.Count; $c++) {
            $tbl1.Cell($r, $c).Shape.TextFrame.TextRange.Font.Bold = $false
        }
    }
}

# ---------------------------------------------------------------
# 2) Investment Summary table: un-bold header row
# ---------------------------------------------------------------
$slideInvestment = Get-SlideByTitle $p "Investment Summary"
$tblInvestment = Get-TableShape $slideInvestment
$tbl2 = $tblInvestment.Table

for ($c = 1; $c -le $tbl2.Columns.Count; $c++) {
    $tbl2.Cell(1, $c).Shape.TextFrame.TextRange.Font.Bold = $false
}

# ---------------------------------------------------------------
# 3) Next Steps: bold the "Label:" prefix of each bullet line
# ---------------------------------------------------------------
$slideNext = Get-SlideByTitle $p "Next Steps"

$bodyShape = $null
for ($j = 1; $j -le $slideNext.Shapes.Count; $j++) {
    $shp = $slideNext.Shapes.Item($j)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "Decision:*") {
        $bodyShape = $shp
    }
}

$tr = $bodyShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs(1, $tr.Text.Length).Count

for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $text = $para.Text
    $colonIdx = $text.IndexOf(":")
    if ($colonIdx -ge 0) {
        $labelLen = $colonIdx + 1
        $labelRange = $para.Characters(1, $labelLen)
        $labelRange.Font.Bold = $true
    }
}
